$wb = $excel.ActiveWorkbook

# Sheet "Overview": Latest HO Xliff Generate Date for 54ef4121-... row (row 3, column G)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-31 10:50:30"

# Sheet "zh-cn": Correspond Handoff/Handback Datetime for 54ef4121-... row (row 3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-31 10:50:26"
$wsZhCn.Range("K3").Value = "2016-08-31 10:50:45"

# Sheet "de-de": Correspond Handoff/Handback Datetime for 54ef4121-... row (row 3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-31 10:50:30"
$wsDeDe.Range("K3").Value = "2016-08-31 10:50:52"
